$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used in the data table (H and N are intentionally skipped/empty)
$cols = @("B","C","D","E","F","I","J","K","L","M")

# New values per row (rows 2-25), in the order of $cols
$data = @{
    2 = @(1.02, 1.0351433820255, 1.03813277289284, 1.033982385555171, 1.043601681888111, 1.035898460803204, 1.040258767711314, 1.040921676020302, 1.036783169240924, 1.046375086769898)
    3 = @(1.02, 1.036899391494863, 1.039492005783241, 1.035522176796944, 1.045435936546052, 1.036417757860136, 1.041654133013537, 1.042089408047707, 1.038130106514441, 1.048017737731556)
    4 = @(1.02, 1.038031744876469, 1.04036810453869, 1.036514585823094, 1.0466194242541, 1.036750536247008, 1.042552890941932, 1.04284107270681, 1.038997249882364, 1.049076826862692)
    5 = @(1.02, 1.038506870827355, 1.040735612094549, 1.036930868047741, 1.047116167797697, 1.036889666701384, 1.042929753695051, 1.043156142837921, 1.039360755311227, 1.049521171726077)
    6 = @(1.02, 1.038586593426845, 1.04079727144843, 1.037000709807369, 1.047199527163981, 1.036912982374282, 1.042992973870695, 1.043208990345301, 1.039421728832572, 1.049595727129627)
    7 = @(1.02, 1.038038097104101, 1.040373018334957, 1.03652015182481, 1.046626064862761, 1.036752398330054, 1.042557930411111, 1.04284528632519, 1.039002111131765, 1.049082767718427)
    8 = @(1.02, 1.035737652975054, 1.038592847383409, 1.034503590730967, 1.044222291686268, 1.036074634761056, 1.040731203718724, 1.041317139139896, 1.037239296460396, 1.046931028937029)
    9 = @(1.02, 1.031653242701108, 1.035429192271198, 1.030919242870373, 1.039959668018877, 1.0348552105366, 1.037479906027776, 1.038593633222114, 1.034098484299864, 1.04310938272528)
    10 = @(1.02, 1.028908446275614, 1.033301239726247, 1.028507851288904, 1.037098664444674, 1.034024961030667, 1.035289638786458, 1.036756508445603, 1.031980433748074, 1.040540342207646)
    11 = @(1.02, 1.027714467743329, 1.032375150395713, 1.027458280989573, 1.035854984633389, 1.03366125773415, 1.034335618579826, 1.035955744095016, 1.031057341669089, 1.039422625868096)
    12 = @(1.02, 1.027270126935029, 1.032030441038019, 1.027067588962825, 1.035392275335609, 1.033525523069556, 1.033980389742595, 1.035657496034355, 1.030713549438906, 1.039006637836861)
    13 = @(1.02, 1.027365478115906, 1.032104415151008, 1.027151431838195, 1.035491562420556, 1.033554667666212, 1.034056626904436, 1.035721508097532, 1.030787335740194, 1.039095906063795)
    14 = @(1.02, 1.027677755746474, 1.032346671367997, 1.027426003414687, 1.035816752398502, 1.03365005094846, 1.034306272973778, 1.035931107402535, 1.031028942467034, 1.039388256988108)
    15 = @(1.02, 1.027870047761477, 1.032495837702088, 1.027595064782062, 1.036017012432532, 1.033708734809939, 1.034459973196715, 1.036060140841229, 1.031177682562599, 1.039568274947285)
    16 = @(1.02, 1.02898756961603, 1.033362601378586, 1.028577391853306, 1.037181099303228, 1.034049009647655, 1.035352833949888, 1.036809540019815, 1.032041569098321, 1.040614407771099)
    17 = @(1.02, 1.029687082532649, 1.033905036865193, 1.029192113836581, 1.037909986499784, 1.034261324983687, 1.03591138400078, 1.037278194464188, 1.032581852728081, 1.041269183963482)
    18 = @(1.02, 1.030094571008381, 1.034220980679786, 1.029550149026859, 1.038334668041715, 1.0343847602621, 1.03623663583355, 1.037551044570244, 1.032896417055785, 1.041650593459926)
    19 = @(1.02, 1.030233425598496, 1.034328633591617, 1.029672141888436, 1.038479395121428, 1.034426780107138, 1.036347447114953, 1.037643993693669, 1.033003578577144, 1.041780558316517)
    20 = @(1.02, 1.029612085946402, 1.033846885298142, 1.029126214098114, 1.037831832171823, 1.0342385874722, 1.035851512956866, 1.037227964946053, 1.032523944913236, 1.041198985620008)
    21 = @(1.02, 1.027585821254375, 1.032275352920488, 1.027345172140365, 1.035721012980363, 1.033621980656174, 1.034232782411811, 1.035869408077715, 1.030957820712593, 1.039302189693359)
    22 = @(1.02, 1.026306933272759, 1.031283102432817, 1.026220519720394, 1.034389499042845, 1.033230593940127, 1.033210018190493, 1.035010543635081, 1.029967834329488, 1.038104852879948)
    23 = @(1.02, 1.026985367794485, 1.031809513673017, 1.026817184939047, 1.035095780677984, 1.03343842904051, 1.033752685564648, 1.035466293382779, 1.030493153806663, 1.038740040867492)
    24 = @(1.02, 1.029645975298942, 1.033873162877171, 1.029155992964968, 1.037867148187955, 1.034248862825131, 1.035878567775574, 1.037250663072848, 1.032550112735214, 1.04123070680629)
    25 = @(1.02, 1.032712924612769, 1.036250332009258, 1.031849652369432, 1.041064963985421, 1.035173479614524, 1.038324379774167, 1.039301446699599, 1.034914648445889, 1.04410103839563)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
